$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "sitp"
$ws.Range("C2").Value = "scooter"
$ws.Range("D2").Value = "motorbike"
$ws.Range("E2").Value = "car"
$ws.Range("F2").Value = "bicycle"
$ws.Range("G2").Value = "2-axis"
$ws.Range("H2").Value = "bus"

$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 27
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 36
$ws.Range("E4").Value = 99
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 1

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 23
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1

$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

$ws.Range("G7").Formula = "=SUM(G3:G6)"
$ws.Range("H7").Formula = "=SUM(H3:H6)"
